$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Insert a new column at A, shifting sw/desc/status/device from A:D to B:E.
$ws.Range("A1").EntireColumn.Insert()

# New "device" column (E) values for the existing rows - order chosen so new
# shared strings are interned as wap(15), dvr(16), atm(17).
$ws.Range("E2").Value = "wap"
$ws.Range("E3").Value = "printer"
$ws.Range("E4").Value = "wap"
$ws.Range("E5").Value = "dvr"
$ws.Range("E6").Value = "atm"

# New row describing the second port on switch 401.
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 401
$ws.Range("C7").Value = "g1/0/2"
$ws.Range("D7").Value = "not"
$ws.Range("E7").Value = "atm"

# New "floor" column (A) header + values, interned last so "floor" becomes 18.
$ws.Range("A1").Value = "floor"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 2

$ws.Range("A1:E7").SetPhonetic()

$ws.Range("D12").Select()
